# Applies the commit's changes:
# - Sheet "BuscaHome": B2 now points to the Tablets category page instead
#   of the Headset product page (both the cell text and its hyperlink),
#   and B2 becomes the selected cell on that sheet.
# - Sheet "Contas": A2 value changes from "MatheusVieira" to "kirxhzees7".
# - Sheet "BuscaLupa": the selected cell changes to A23 (view only).
# - The sheet that was active before the edit ("Contas") stays active.

$wb = $excel.ActiveWorkbook

# --- Sheet "BuscaHome" ---
# Update the shared string used by B2 first, so the new shared string
# (the Tablets URL) reuses the slot vacated by the old "MatheusVieira"
# text before the brand-new "kirxhzees7" string gets appended below.
$wsBuscaHome = $wb.Worksheets.Item("BuscaHome")
$wsBuscaHome.Range("B2").Value = "https://www.advantageonlineshopping.com/#/category/Tablets/3"

# Rebuild both hyperlinks on this sheet (B2 and A2) so B2 points at the
# new Tablets category instead of the old Headset product page; A2's
# hyperlink target is unchanged.
$wsBuscaHome.Hyperlinks.Delete()
$wsBuscaHome.Hyperlinks.Add($wsBuscaHome.Range("B2"), "https://www.advantageonlineshopping.com/", "/category/Tablets/3")
$wsBuscaHome.Hyperlinks.Add($wsBuscaHome.Range("A2"), "https://www.advantageonlineshopping.com/", "/product/18")

# Restore the "Hiperlink" style on both cells (Hyperlinks.Add resets it).
$wsBuscaHome.Range("A2").Style = "Hiperlink"
$wsBuscaHome.Range("B2").Style = "Hiperlink"

$wsBuscaHome.Range("B2").Select()

# --- Sheet "Contas" ---
$wsContas = $wb.Worksheets.Item("Contas")
$wsContas.Range("A2").Value = "kirxhzees7"

# --- Sheet "BuscaLupa" (selection only) ---
$wsBuscaLupa = $wb.Worksheets.Item("BuscaLupa")
$wsBuscaLupa.Range("A23").Select()

# Restore the originally active sheet/tab.
$wsContas.Activate()
